$d = $word.ActiveDocument

# --- 1. Create the three new character styles -----------------------------

$ganStyle = $d.Styles.Add("GaNStyle", 2)       # wdStyleTypeCharacter
$ganStyle.Font.Name = "Calibri"
$ganStyle.Font.Size = 14

$ganParagraph = $d.Styles.Add("GaNParagraph", 2)
$ganParagraph.Font.Name = "Calibri"
$ganParagraph.Font.Size = 10

$ganLinks = $d.Styles.Add("GaNLinks", 2)
$ganLinks.Font.Name = "Calibri"
$ganLinks.Font.Size = 9.5
$ganLinks.Font.Bold = $true
$ganLinks.Font.Underline = 1          # wdUnderlineSingle
$ganLinks.Font.Color = 8388608        # wdColorNavy -> w:color val="000080"

# --- 2. Apply the styles to the matching runs ------------------------------

function Apply-StyleToAllMatches($searchText, $styleName) {
    $rng = $d.Content
    $rng.Start = 0
    while ($rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
        $rng.Style = $styleName
        $rng.Collapse(0)
    }
}

Apply-StyleToAllMatches "V roku 2022 môžete pozorovať Súhvezdie Bootes: 14. – 23. mája, 13. – 22. júna, 12. – 21. júla" "GaNStyle"

Apply-StyleToAllMatches "Stávate sa súčasťou celosvetovej kampane Globe at Night, ktorej cieľom je meranie svetelného znečistenia. Pozorovaním  Súhvezdie Bootes na nočnej oblohe a porovnávaním skutočnej situácie s našimi mapkami sa nielenže dozviete, ako osvetlenie vo Vašom okolí prispieva k svetelnému znečisteniu, ale budete môcť porovnať úroveň svetelného znečistenia aj s inými lokalitami z celého sveta. Vaše pozorovanie tiež rozšíri online databázu dokumentujúcu viditeľnosť nočnej oblohy na našej planéte" "GaNParagraph"

Apply-StyleToAllMatches "Mapky v tomto dokumente pripravil Jan Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)." "GaNLinks"
